$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 1940
$ws.Range("L3").Value = 1969
$ws.Range("L4").Value = 547
$ws.Range("L5").Value = 114
$ws.Range("L6").Value = 1782
$ws.Range("L7").Value = 6352

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L4").Value = 22
$ws.Range("L5").Value = 21
$ws.Range("L7").Value = 214
$ws.Range("L8").Value = 395
$ws.Range("L12").Value = 15
$ws.Range("L14").Value = 31
$ws.Range("L19").Value = 180
$ws.Range("L20").Value = 167
$ws.Range("L21").Value = 21
$ws.Range("L25").Value = 34
$ws.Range("L29").Value = 327
$ws.Range("L31").Value = 63
$ws.Range("L33").Value = 281
$ws.Range("L34").Value = 41
$ws.Range("L35").Value = 12
$ws.Range("L36").Value = 93
$ws.Range("L37").Value = 228
$ws.Range("L42").Value = 199
$ws.Range("L43").Value = 51
$ws.Range("L47").Value = 42
$ws.Range("L49").Value = 34
$ws.Range("L50").Value = 39
$ws.Range("L51").Value = 75
$ws.Range("L52").Value = 130
$ws.Range("L53").Value = 76
$ws.Range("L54").Value = 135
$ws.Range("L55").Value = 58
$ws.Range("L57").Value = 29
$ws.Range("L58").Value = 6
$ws.Range("L60").Value = 35
$ws.Range("L63").Value = 19
$ws.Range("L64").Value = 47
$ws.Range("L65").Value = 119
$ws.Range("L67").Value = 224
$ws.Range("L68").Value = 17
$ws.Range("L69").Value = 13
$ws.Range("L72").Value = 28
$ws.Range("L73").Value = 52
$ws.Range("L75").Value = 26
$ws.Range("L76").Value = 66
$ws.Range("L78").Value = 87
$ws.Range("L79").Value = 177
$ws.Range("L83").Value = 152
$ws.Range("L84").Value = 64
$ws.Range("L85").Value = 339
$ws.Range("L86").Value = 49
$ws.Range("L89").Value = 77
$ws.Range("L90").Value = 61
$ws.Range("L91").Value = 87
$ws.Range("L94").Value = 77
$ws.Range("L95").Value = 94
$ws.Range("L99").Value = 99
$ws.Range("L101").Value = 6352

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 62
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 214

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 105
$ws.Range("L3").Value = 141
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 339

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 13

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 10
$ws.Range("L7").Value = 76

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 110
$ws.Range("L3").Value = 134
$ws.Range("L6").Value = 107
$ws.Range("L7").Value = 395

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L6").Value = 36
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L6").Value = 97
$ws.Range("L7").Value = 281

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 32
$ws.Range("L7").Value = 94

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 73
$ws.Range("L3").Value = 62
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 228

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L3").Value = 37
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L2").Value = 27
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 72
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 224

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("L2").Value = 25
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L2").Value = 8
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 105
$ws.Range("L3").Value = 116
$ws.Range("L7").Value = 327

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 57
$ws.Range("L3").Value = 55
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 180

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 50
$ws.Range("L3").Value = 55
$ws.Range("L4").Value = 21
$ws.Range("L6").Value = 68
$ws.Range("L7").Value = 199

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 28
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L2").Value = 24
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L2").Value = 16
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 47

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L2").Value = 18
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item('Gold Coast')
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("L3").Value = 14
$ws.Range("L4").Value = 5
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 21

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 28
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("L2").Value = 15
$ws.Range("L7").Value = 26

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("L2").Value = 13
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 15

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 6
